# Update the acquisition timestamp (column A) for rows 2-16 on the
# active sheet ("ランサーズ") from 2025-11-10 01:23:14 to 2025-11-10 01:54:22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-10 01:54:22"

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
